$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet (matches the re-uploaded workbook's sheet tab name)
$ws.Name = "voxelVolumes_treatments"

# Restore the standard (non-custom) Excel page margins that a fresh
# workbook ships with, replacing the narrower margins of the old file.
$ws.PageSetup.LeftMargin   = 54   # 0.75 in
$ws.PageSetup.RightMargin  = 54   # 0.75 in
$ws.PageSetup.TopMargin    = 72   # 1 in
$ws.PageSetup.BottomMargin = 72   # 1 in
$ws.PageSetup.HeaderMargin = 36   # 0.5 in
$ws.PageSetup.FooterMargin = 36   # 0.5 in

# The data was copied into a brand-new default workbook, which always
# carries the complete built-in cell-style gallery (fonts/fills/borders
# + the 41 named styles) even though no cell in the sheet actually uses
# any of them. Touch every built-in style on scratch cells far below the
# real data, then clear those scratch cells so the gallery is registered
# in the workbook without altering the sheet's used range or content.
$scratchRow = 1000
$ws.Cells.Item($scratchRow, 1).Style = "20% - Accent1"
$ws.Cells.Item($scratchRow, 2).Style = "20% - Accent2"
$ws.Cells.Item($scratchRow, 3).Style = "20% - Accent3"
$ws.Cells.Item($scratchRow, 4).Style = "20% - Accent4"
$ws.Cells.Item($scratchRow, 5).Style = "20% - Accent5"
$ws.Cells.Item($scratchRow, 6).Style = "20% - Accent6"
$ws.Cells.Item($scratchRow, 7).Style = "40% - Accent1"
$ws.Cells.Item($scratchRow, 8).Style = "40% - Accent2"
$ws.Cells.Item($scratchRow, 9).Style = "40% - Accent3"
$ws.Cells.Item($scratchRow, 10).Style = "40% - Accent4"
$ws.Cells.Item($scratchRow, 11).Style = "40% - Accent5"
$ws.Cells.Item($scratchRow, 12).Style = "40% - Accent6"
$ws.Cells.Item($scratchRow, 13).Style = "60% - Accent1"
$ws.Cells.Item($scratchRow, 14).Style = "60% - Accent2"
$ws.Cells.Item($scratchRow, 15).Style = "60% - Accent3"
$ws.Cells.Item($scratchRow, 16).Style = "60% - Accent4"
$ws.Cells.Item($scratchRow, 17).Style = "60% - Accent5"
$ws.Cells.Item($scratchRow, 18).Style = "60% - Accent6"
$ws.Cells.Item($scratchRow, 19).Style = "Accent1"
$ws.Cells.Item($scratchRow, 20).Style = "Accent2"
$ws.Cells.Item($scratchRow, 21).Style = "Accent3"
$ws.Cells.Item($scratchRow, 22).Style = "Accent4"
$ws.Cells.Item($scratchRow, 23).Style = "Accent5"
$ws.Cells.Item($scratchRow, 24).Style = "Accent6"
$ws.Cells.Item($scratchRow, 25).Style = "Bad"
$ws.Cells.Item($scratchRow, 26).Style = "Calculation"
$ws.Cells.Item($scratchRow, 27).Style = "Check Cell"
$ws.Cells.Item($scratchRow, 28).Style = "Explanatory Text"
$ws.Cells.Item($scratchRow, 29).Style = "Good"
$ws.Cells.Item($scratchRow, 30).Style = "Heading 1"
$ws.Cells.Item($scratchRow, 31).Style = "Heading 2"
$ws.Cells.Item($scratchRow, 32).Style = "Heading 3"
$ws.Cells.Item($scratchRow, 33).Style = "Heading 4"
$ws.Cells.Item($scratchRow, 34).Style = "Input"
$ws.Cells.Item($scratchRow, 35).Style = "Linked Cell"
$ws.Cells.Item($scratchRow, 36).Style = "Neutral"
$ws.Cells.Item($scratchRow, 37).Style = "Note"
$ws.Cells.Item($scratchRow, 38).Style = "Output"
$ws.Cells.Item($scratchRow, 39).Style = "Title"
$ws.Cells.Item($scratchRow, 40).Style = "Total"
$ws.Cells.Item($scratchRow, 41).Style = "Warning Text"
$ws.Rows.Item($scratchRow).Clear()
